# chore: Javadoc for TestMain
#
# TrainSearchData.xlsx test-data tweak:
#  - "SearchData" sheet: the " Date" column becomes an "After Days" column
#    (a plain day-offset number instead of a live TODAY()+N formula).
#  - "Testing" sheet: the logged "Date" row's Expected/Actual formula now
#    reads that offset from SearchData!C2 instead of a hard-coded "+4".

$wb = $excel.ActiveWorkbook

$search = $wb.Worksheets.Item("SearchData")
$testing = $wb.Worksheets.Item("Testing")

# --- SearchData sheet -------------------------------------------------
# Header C1: " Date" -> "After Days"
$search.Range("C1").Value = "After Days"

# C2: was a live formula TEXT(TODAY()+4, "dd-mmm-yy ddd"); now a plain
# number of days to add ("After Days"). Clear the old date number format
# so the cell reads as a plain integer, not a date serial.
$search.Range("C2").Style = "Normal"
$search.Range("C2").Value = 6

# Selection moves to the edited cell.
$search.Range("C2").Select()

# --- Testing sheet ------------------------------------------------------
# B6 formula now offsets TODAY() by SearchData!C2 instead of a literal 4.
$testing.Range("B6").Formula = '=TEXT(TODAY()+SearchData!C2, "dd-mmm-yy ddd")'

# C6 ("Actual") is a plain logged copy of the freshly computed B6 value,
# mirroring how C2:C5 mirror B2:B5 on this sheet.
$testing.Range("C6").Value = $testing.Range("B6").Value2

# Selection moves to the edited cell.
$testing.Range("B6").Select()
